# Fruta / hortaliza, semanal
# Insert a new data row at row 228 (pushing the existing rows 228-334 down to
# 229-335) and populate the newly inserted row with a new price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 228..334 down to 229..335, leaving a blank row 228 behind.
$ws.Rows.Item(228).Insert()

# Fill in the new row 228 with the new observation.
$ws.Range("A228").Value = 4
$ws.Range("B228").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C228").Value = "Los Lagos"
$ws.Range("D228").Value = 44609
$ws.Range("E228").Value = 10
$ws.Range("F228").Value = 100114001
$ws.Range("G228").Value = "Papa"
$ws.Range("H228").Value = "Patagonia"
$ws.Range("I228").Value = "1a nueva(o)"
$ws.Range("J228").Value = 250
$ws.Range("K228").Value = 6000
$ws.Range("L228").Value = 7000
$ws.Range("M228").Value = 6400
$ws.Range("N228").Value = "`$/saco 25 kilos"
$ws.Range("O228").Value = "Provincia de Llanquihue"
$ws.Range("P228").Value = 256
$ws.Range("Q228").Value = 25
$ws.Range("R228").Value = "Hortaliza"
